$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 878.3333
$ws.Range("I12").Value = 354
$ws.Range("K12").Value = 354
$ws.Range("M12").Value = -184

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1440.3636
$ws.Range("I39").Value = 65
$ws.Range("J39").Value = 2586.5
$ws.Range("K39").Value = 195
$ws.Range("L39").Value = 7759.5
$ws.Range("M39").Value = 101
$ws.Range("N39").Value = -8351.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 9641.429
$ws.Range("I74").Value = 9122.5
$ws.Range("K74").Value = 9122.5
$ws.Range("M74").Value = -8186.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 9641.429
$ws.Range("I77").Value = 9122.5
$ws.Range("K77").Value = 45612.5
$ws.Range("M77").Value = -40932.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 79334.30499999999
$ws.Range("J112").Value = 102807.1
$ws.Range("L112").Value = 308421.3
$ws.Range("N112").Value = -310637.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3509.4211
$ws.Range("J125").Value = 5025.8184
$ws.Range("L125").Value = 45232.3656
$ws.Range("N125").Value = -50152.3656

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3556.93
$ws.Range("J138").Value = 3682.4363
$ws.Range("L138").Value = 11047.3089
$ws.Range("N138").Value = -21327.3089

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1842.4717
$ws.Range("I132").Value = 1635.9348
$ws.Range("K132").Value = 4907.8044
$ws.Range("M132").Value = -2377.8044

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2850.875
$ws.Range("I20").Value = 2831.2
$ws.Range("K20").Value = 2831.2
$ws.Range("M20").Value = -2584.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1994.7931
$ws.Range("I105").Value = 1709.7273
$ws.Range("J105").Value = 2890.7144
$ws.Range("K105").Value = 1709.7273
$ws.Range("L105").Value = 2890.7144
$ws.Range("M105").Value = 37.27269999999999
$ws.Range("N105").Value = -6384.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2669465.8
$ws.Range("I134").Value = 3335332.8
$ws.Range("K134").Value = 10005998.4
$ws.Range("M134").Value = -10003463.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1461.8572
$ws.Range("J31").Value = 1551.898
$ws.Range("L31").Value = 1551.898
$ws.Range("N31").Value = -2141.898

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1461.8572
$ws.Range("J34").Value = 1551.898
$ws.Range("L34").Value = 1551.898
$ws.Range("N34").Value = -1955.898

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("N54").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1727.4615
$ws.Range("I68").Value = 1398.7
$ws.Range("J68").Value = 1932.9375
$ws.Range("K68").Value = 4196.1
$ws.Range("L68").Value = 5798.8125
$ws.Range("M68").Value = -3385.1
$ws.Range("N68").Value = -7420.8125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1727.4615
$ws.Range("I71").Value = 1398.7
$ws.Range("J71").Value = 1932.9375
$ws.Range("K71").Value = 12588.3
$ws.Range("L71").Value = 17396.4375
$ws.Range("M71").Value = -8532.300000000001
$ws.Range("N71").Value = -25508.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 225.21428
$ws.Range("I86").Value = 216
$ws.Range("J86").Value = 241.8
$ws.Range("K86").Value = 648
$ws.Range("L86").Value = 725.4000000000001
$ws.Range("M86").Value = 538
$ws.Range("N86").Value = -3097.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 225.21428
$ws.Range("I89").Value = 216
$ws.Range("J89").Value = 241.8
$ws.Range("K89").Value = 1944
$ws.Range("L89").Value = 2176.2
$ws.Range("M89").Value = 3984
$ws.Range("N89").Value = -14032.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 4650
$ws.Range("I95").Value = 1800
$ws.Range("J95").Value = 7500
$ws.Range("K95").Value = 5400
$ws.Range("L95").Value = 22500
$ws.Range("M95").Value = -3341
$ws.Range("N95").Value = -26618

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 726.6667
$ws.Range("I116").Value = 726.6667
$ws.Range("K116").Value = 2180.0001
$ws.Range("M116").Value = 1261.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 95.5
$ws.Range("I118").Value = 95.5
$ws.Range("K118").Value = 286.5
$ws.Range("M118").Value = 956.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 15284.143
$ws.Range("I121").Value = 230
$ws.Range("J121").Value = 26574.75
$ws.Range("K121").Value = 690
$ws.Range("L121").Value = 79724.25
$ws.Range("M121").Value = 620
$ws.Range("N121").Value = -82344.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 992.3333
$ws.Range("I97").Value = 821.1923
$ws.Range("J97").Value = 2104.75
$ws.Range("K97").Value = 821.1923
$ws.Range("L97").Value = 2104.75
$ws.Range("M97").Value = -325.1923
$ws.Range("N97").Value = -3096.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1410.5
$ws.Range("J102").Value = 1494.8334
$ws.Range("L102").Value = 1494.8334
$ws.Range("N102").Value = -4738.8334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2628.913
$ws.Range("I122").Value = 2879.25
$ws.Range("K122").Value = 8637.75
$ws.Range("M122").Value = -6187.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3070.8372
$ws.Range("I132").Value = 2487.3794
$ws.Range("J132").Value = 4279.4287
$ws.Range("K132").Value = 7462.138199999999
$ws.Range("L132").Value = 12838.2861
$ws.Range("M132").Value = -4932.138199999999
$ws.Range("N132").Value = -17898.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2761.647
$ws.Range("I22").Value = 2489.9
$ws.Range("J22").Value = 3149.8572
$ws.Range("K22").Value = 2489.9
$ws.Range("L22").Value = 3149.8572
$ws.Range("M22").Value = -2194.9
$ws.Range("N22").Value = -3739.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2761.647
$ws.Range("I27").Value = 2489.9
$ws.Range("J27").Value = 3149.8572
$ws.Range("K27").Value = 2489.9
$ws.Range("L27").Value = 3149.8572
$ws.Range("M27").Value = -2382.9
$ws.Range("N27").Value = -3363.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2933.28
$ws.Range("J46").Value = 3064.652
$ws.Range("L46").Value = 3064.652
$ws.Range("N46").Value = -3440.652

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 806.6896400000001
$ws.Range("I55").Value = 954.6667
$ws.Range("J55").Value = 564.5454999999999
$ws.Range("K55").Value = 954.6667
$ws.Range("L55").Value = 564.5454999999999
$ws.Range("M55").Value = -781.6667
$ws.Range("N55").Value = -910.5454999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8544.4
$ws.Range("I132").Value = 9493.125
$ws.Range("K132").Value = 28479.375
$ws.Range("M132").Value = -25949.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 94999.5
$ws.Range("J133").Value = 94999.5
$ws.Range("L133").Value = 94999.5
$ws.Range("N133").Value = -100059.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9823.380999999999
$ws.Range("J81").Value = 9639.4
$ws.Range("L81").Value = 19278.8
$ws.Range("N81").Value = -21400.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 9823.380999999999
$ws.Range("J84").Value = 9639.4
$ws.Range("L84").Value = 96394
$ws.Range("N84").Value = -107002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2917.45
$ws.Range("I126").Value = 2828.0908
$ws.Range("J126").Value = 3338.7144
$ws.Range("K126").Value = 8484.2724
$ws.Range("L126").Value = 10016.1432
$ws.Range("M126").Value = -6014.2724
$ws.Range("N126").Value = -14956.1432
